# Generate Report for Handoff
# Renames the existing handed-off file (a52f26a9-...) to d7f33c70-...,
# refreshes its timestamps/xlf hashes, and adds a brand new duplicate
# record (ffff053ed6fe-...) as row 3 on every sheet.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/59d4203842c326fcf2a2a73d4199acfeada3f6c2/e2e/"

$oldName = "a52f26a9-04f3-4155-93ce-0b26dd45f53e"
$newName = "d7f33c70-05aa-4ad7-a612-d574772ec7d8"
$dupName = "ffff053ed6fe-9668-4ed5-abf2-084c2d724b4d"

$oldHash = "8618f1a63b792b4efe4c19a2fd11272ef5d05485"
$newHash = "b21d747795926e5e9f8f28898e3a20c811cdb8c9"

$dateOverviewZh = "2017-02-09 09:11:32"
$dateZhXlf       = "2017-02-09 09:11:00"
$dateDeXlf       = "2017-02-09 09:11:32"

# ---------------------------------------------------------------
# Sheet "Overview" (sheet1)
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newName.md"
$wsOverview.Range("G2").Value = $dateOverviewZh

$wsOverview.Range("A3").Value = "$dupName.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $dateOverviewZh
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$baseUrl$newName.md", [Type]::Missing, [Type]::Missing, "e2e\$newName.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$baseUrl$dupName.md", [Type]::Missing, [Type]::Missing, "e2e\$dupName.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$newZhXlf = "$newName.$newHash.zh-cn.xlf"

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $dateZhXlf
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("F3").Style = "Normal"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $dateZhXlf
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("O3").Value = "'True"
$wsZh.Range("O3").Style = "Normal"
$wsZh.Range("Q3").Value = "'False"
$wsZh.Range("Q3").Style = "Normal"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$baseUrl$newName.md", [Type]::Missing, [Type]::Missing, "$newName.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$baseUrl$dupName.md", [Type]::Missing, [Type]::Missing, "$dupName.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:R3"))

# ---------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$newDeXlf = "$newName.$newHash.de-de.xlf"

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $dateDeXlf
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("F3").Style = "Normal"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $dateDeXlf
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("O3").Value = "'True"
$wsDe.Range("O3").Style = "Normal"
$wsDe.Range("Q3").Value = "'False"
$wsDe.Range("Q3").Style = "Normal"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$baseUrl$newName.md", [Type]::Missing, [Type]::Missing, "$newName.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$baseUrl$dupName.md", [Type]::Missing, [Type]::Missing, "$dupName.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:R3"))

Write-Host "Done"
